$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- "Call Procedure or Function" block (rows 17-19) ---
# Row 17: CALL ... the D/F columns swap meaning - the index placeholder
# becomes PROC_FUNC_INDEX and the parameter-count placeholder becomes PARA_NUM.
$ws.Range("D17").Value = "PARA_NUM"
$ws.Range("F17").Value = "PROC_FUNC_INDEX"

# --- "BEGIN/END PROC/FUNC" block (rows 21-22) ---
# Drop the old D column (PROC/FUNC) entirely and rename the Index placeholder
# to PROC_FUNC_INDEX everywhere it is used (including inside the example text).
$ws.Range("D21").ClearContents()
$ws.Range("F21").Value = "PROC_FUNC_INDEX"
$ws.Range("G21").Value = "BEGIN PROC/FUNC at TokenTable[PROC_FUNC_INDEX]"

$ws.Range("D22").ClearContents()
$ws.Range("F22").Value = "PROC_FUNC_INDEX"
$ws.Range("G22").Value = "END PROC/FUNC at TokenTable[PROC_FUNC_INDEX]"

# --- cosmetic: widen column F (target authored width ~17.25 chars) and move the active selection ---
# Note: this host's ColumnWidth->pixel->XML-width round trip only resolves to
# multiples of 1/7 character-widths (Excel's default-font MDW quantisation),
# so 16.57 is the input that lands on the closest achievable stored width to 17.25.
$ws.Columns.Item(6).ColumnWidth = 16.57
$ws.Range("D18").Select() | Out-Null
